$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = -0.17212692718101152
$ws.Cells.Item(1, 2).Value2 = 0.17178198059662719
$ws.Cells.Item(2, 1).Value2 = -0.1141868407097153
$ws.Cells.Item(2, 2).Value2 = 0.11339668917920509
$ws.Cells.Item(3, 1).Value2 = -0.063684109576040981
$ws.Cells.Item(3, 2).Value2 = 0.063544102191725571
$ws.Cells.Item(4, 1).Value2 = -0.055544102234188486
$ws.Cells.Item(4, 2).Value2 = 0.055434557598919909
$ws.Cells.Item(5, 1).Value2 = -0.099412236405023613
$ws.Cells.Item(5, 2).Value2 = 0.098015769466182512
$ws.Cells.Item(6, 1).Value2 = -0.052695471012532025
$ws.Cells.Item(6, 2).Value2 = 0.052093980940645324
$ws.Cells.Item(7, 1).Value2 = -0.042093981000257408
$ws.Cells.Item(7, 2).Value2 = 0.041942291428740397
$ws.Cells.Item(8, 1).Value2 = -0.03194229149052763
$ws.Cells.Item(8, 2).Value2 = 0.031657157480617482
$ws.Cells.Item(9, 1).Value2 = -0.029657157508625076
$ws.Cells.Item(9, 2).Value2 = 0.02941759997664084
$ws.Cells.Item(10, 1).Value2 = -0.027417600006662823
$ws.Cells.Item(10, 2).Value2 = 0.027401476909778566
$ws.Cells.Item(11, 1).Value2 = -0.024401476944639455
$ws.Cells.Item(11, 2).Value2 = 0.024374415891500689
$ws.Cells.Item(12, 1).Value2 = -0.020874415929232004
$ws.Cells.Item(12, 2).Value2 = 0.020675250434531112
$ws.Cells.Item(13, 1).Value2 = -0.017175250474331882
$ws.Cells.Item(13, 2).Value2 = 0.017084495788584952
$ws.Cells.Item(14, 1).Value2 = -0.00908449584914095
$ws.Cells.Item(14, 2).Value2 = 0.0090547995530263492
$ws.Cells.Item(15, 1).Value2 = -0.0080547995823536667
$ws.Cells.Item(15, 2).Value2 = 0.0080356205135814207
$ws.Cells.Item(16, 1).Value2 = -0.0060356205477782865
$ws.Cells.Item(16, 2).Value2 = 0.006003625841610738
$ws.Cells.Item(17, 1).Value2 = -0.0040036258761961818
$ws.Cells.Item(17, 2).Value2 = 0.003999999956245226
$ws.Cells.Item(18, 1).Value2 = -0.016105902257415039
$ws.Cells.Item(18, 2).Value2 = 0.016091685542281908
$ws.Cells.Item(19, 1).Value2 = -0.01209168556103446
$ws.Cells.Item(19, 2).Value2 = 0.012016774607611413
$ws.Cells.Item(20, 1).Value2 = -0.0080167746276877949
$ws.Cells.Item(20, 2).Value2 = 0.008005661366770056
$ws.Cells.Item(21, 1).Value2 = -0.0040056613870742552
$ws.Cells.Item(21, 2).Value2 = 0.0039999999795261587
$ws.Cells.Item(22, 1).Value2 = -0.045712579626080441
$ws.Cells.Item(22, 2).Value2 = 0.045499060643219735
$ws.Cells.Item(23, 1).Value2 = -0.040499060672146037
$ws.Cells.Item(23, 2).Value2 = 0.040099049529761288
$ws.Cells.Item(24, 1).Value2 = -0.020099049628972132
$ws.Cells.Item(24, 2).Value2 = 0.019999999899435572
$ws.Cells.Item(25, 1).Value2 = -0.043320298476842822
$ws.Cells.Item(25, 2).Value2 = 0.043304066955631271
$ws.Cells.Item(26, 1).Value2 = -0.070239130639169645
$ws.Cells.Item(26, 2).Value2 = 0.070140240261634901
$ws.Cells.Item(27, 1).Value2 = -0.067640240289138731
$ws.Cells.Item(27, 2).Value2 = 0.067069099611548832
$ws.Cells.Item(28, 1).Value2 = -0.065069099642643735
$ws.Cells.Item(28, 2).Value2 = 0.064689857745878321
$ws.Cells.Item(29, 1).Value2 = -0.057689857803662647
$ws.Cells.Item(29, 2).Value2 = 0.057588916796777134
$ws.Cells.Item(30, 1).Value2 = 0.0024110829046515292
$ws.Cells.Item(30, 2).Value2 = -0.0024667788476895325
$ws.Cells.Item(31, 1).Value2 = -0.049768870180074032
$ws.Cells.Item(31, 2).Value2 = 0.049619232490785592
$ws.Cells.Item(32, 1).Value2 = -0.039619232565467399
$ws.Cells.Item(32, 2).Value2 = 0.039507653281411592

# Set column widths to match target (closest achievable via ColumnWidth rounding)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
